$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$statQuery = @'
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Bisulfite-Seq" IN es
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s:study)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Bisulfite-Seq" IN es
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Bisulfite-Seq" IN es
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`
'@

$participantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (samp)<--(f:file)
WITH p, samp, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Bisulfite-Seq" IN es
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, "Not specified in data"))) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

$sampleQuery = @'
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Bisulfite-Seq" IN es
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(samp.sample_tumor_status,'') as `Tumor`,
    coalesce(samp.sample_type,'') as `Analyte Type`
ORDER BY samp.sample_id LIMIT 100
'@

$fileQuery = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Bisulfite-Seq" IN es
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name LIMIT 100
'@

# StatQuery shared by all three tabs - set first so it gets the lowest new shared-string index
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# Row 2: ParticipantsTab
$ws.Range("B2").Value = $participantQuery

# Row 3: SamplesTab
$ws.Range("B3").Value = $sampleQuery

# Row 4: FilesTab
$ws.Range("B4").Value = $fileQuery

# The longer query text now wraps to the maximum row height Excel allows
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B4").Select()
